# Apply the cryptos-list price/volume/coin updates described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.401.00"
$ws.Range("E2").Value = "  +3.55%  "
$ws.Range("D3").Value = "1.793.89"
$ws.Range("E3").Value = "  +4.06%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "336.70"
$ws.Range("D5").Style = $ws.Range("C5").Style
$ws.Range("E5").Value = "  +1.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.0000"
$ws.Range("D6").Style = $ws.Range("C6").Style
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3795"
$ws.Range("D7").Style = $ws.Range("C7").Style
$ws.Range("E7").Value = "  +1.76%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3431"
$ws.Range("D8").Style = $ws.Range("C8").Style
$ws.Range("E8").Value = "  +1.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.03"
$ws.Range("D9").Style = $ws.Range("C9").Style
$ws.Range("E9").Value = "  -0.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.201"
$ws.Range("D10").Style = $ws.Range("C10").Style
$ws.Range("E10").Value = "  +1.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07496"
$ws.Range("D11").Style = $ws.Range("C11").Style
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("D12").Style = $ws.Range("C12").Style
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.05"
$ws.Range("D13").Style = $ws.Range("C13").Style
$ws.Range("E13").Value = "  +9.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.477"
$ws.Range("D14").Style = $ws.Range("C14").Style
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("D15").Value = "1.792.12"
$ws.Range("E15").Value = "  +4.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.053"
$ws.Range("D16").Style = $ws.Range("C16").Style
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001097"
$ws.Range("D17").Style = $ws.Range("C17").Style
$ws.Range("E17").Value = "  +2.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06652"
$ws.Range("D18").Style = $ws.Range("C18").Style
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "84.69"
$ws.Range("D19").Style = $ws.Range("C19").Style
$ws.Range("E19").Value = "  +3.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.0000"
$ws.Range("D20").Style = $ws.Range("C20").Style
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.37"
$ws.Range("D21").Style = $ws.Range("C21").Style
$ws.Range("E21").Value = "  +5.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.481"
$ws.Range("D22").Style = $ws.Range("C22").Style
$ws.Range("E22").Value = "  +5.19%  "
$ws.Range("D23").Value = "27.378.53"
$ws.Range("E23").Value = "  +3.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.50"
$ws.Range("D24").Style = $ws.Range("C24").Style
$ws.Range("E24").Value = "  -1.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.443"
$ws.Range("D25").Style = $ws.Range("C25").Style
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.580"
$ws.Range("D26").Style = $ws.Range("C26").Style
$ws.Range("E26").Value = "  +8.31%  "
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.496"
$ws.Range("D27").Style = $ws.Range("C27").Style
$ws.Range("E27").Value = "  +6.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.38"
$ws.Range("D28").Style = $ws.Range("C28").Style
$ws.Range("E28").Value = "  +9.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "150.93"
$ws.Range("D29").Style = $ws.Range("C29").Style
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "1.992.11"
$ws.Range("E30").Value = "  +4.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "133.47"
$ws.Range("D31").Style = $ws.Range("C31").Style
$ws.Range("E31").Value = "  +1.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.057"
$ws.Range("D32").Style = $ws.Range("C32").Style
$ws.Range("E32").Value = "  -0.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.129"
$ws.Range("D33").Style = $ws.Range("C33").Style
$ws.Range("E33").Value = "  +2.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08711"
$ws.Range("D34").Style = $ws.Range("C34").Style
$ws.Range("E34").Value = "  +1.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.30"
$ws.Range("D35").Style = $ws.Range("C35").Style
$ws.Range("E35").Value = "  +4.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.676"
$ws.Range("D36").Style = $ws.Range("C36").Style
$ws.Range("E36").Value = "  -0.93%  "
$ws.Range("B37").Value = "TheSandbox"
$ws.Range("C37").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6898"
$ws.Range("D37").Style = $ws.Range("C37").Style
$ws.Range("E37").Value = "  +11.26%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.440"
$ws.Range("D38").Style = $ws.Range("C38").Style
$ws.Range("E38").Value = "  +1.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.865"
$ws.Range("D39").Style = $ws.Range("C39").Style
$ws.Range("E39").Value = "  +5.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06350"
$ws.Range("D40").Style = $ws.Range("C40").Style
$ws.Range("E40").Value = "  +2.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2201"
$ws.Range("D41").Style = $ws.Range("C41").Style
$ws.Range("E41").Value = "  +2.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02338"
$ws.Range("D42").Style = $ws.Range("C42").Style
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("E43").Value = "  +4.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.36"
$ws.Range("D44").Style = $ws.Range("C44").Style
$ws.Range("E44").Value = "  +1.17%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6454"
$ws.Range("D45").Style = $ws.Range("C45").Style
$ws.Range("E45").Value = "  +7.39%  "
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9998"
$ws.Range("D46").Style = $ws.Range("C46").Style
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.855"
$ws.Range("D47").Style = $ws.Range("C47").Style
$ws.Range("E47").Value = "  -1.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.117"
$ws.Range("D48").Style = $ws.Range("C48").Style
$ws.Range("E48").Value = "  +3.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "129.94"
$ws.Range("D49").Style = $ws.Range("C49").Style
$ws.Range("E49").Value = "  +1.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07182"
$ws.Range("D50").Style = $ws.Range("C50").Style
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.16"
$ws.Range("D51").Style = $ws.Range("C51").Style
$ws.Range("E51").Value = "  +3.05%  "
